# Clear the "My Marks" entries (D5:D10) on Sheet1 - the dependent formulas in
# D11:D13 (Total / Letter Grade / Grade Point) will recalc automatically once
# the inputs are cleared.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5:D10").ClearContents()

# Move the active selection from F10 to F9 as recorded in the saved view.
$ws.Range("F9").Select()
